# Updates cryptos list price (D) and 1h volume/change (E) columns
# per the Thu May  2 08:50:23 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="57.944.87"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").Formula = '="2.946.20"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = "  +3.07%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Formula = '="552.83"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +1.61%  "

$ws.Range("D6").Formula = '="133.18"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +10.77%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  +5.59%  "

$ws.Range("D9").Formula = '="2.937.24"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +3.09%  "

$ws.Range("E10").Value = "  +3.80%  "

$ws.Range("E11").Value = "  +1.91%  "

$ws.Range("E12").Value = "  +5.69%  "

$ws.Range("E13").Value = "  +6.02%  "

$ws.Range("D14").Formula = '="32.77"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +6.19%  "

$ws.Range("E15").Value = "  +3.82%  "

$ws.Range("D16").Formula = '="3.429.09"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +3.17%  "

$ws.Range("D17").Formula = '="6.88"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  +10.74%  "

$ws.Range("D18").Formula = '="2.937.22"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +2.91%  "

$ws.Range("D19").Formula = '="57.915.44"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +1.55%  "

$ws.Range("D20").Formula = '="417.20"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +2.65%  "

$ws.Range("D21").Formula = '="13.32"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +6.06%  "

$ws.Range("D22").Formula = '="0.694"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +8.63%  "

$ws.Range("D23").Formula = '="13.40"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +8.89%  "

$ws.Range("D24").Formula = '="7.00"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +4.91%  "

$ws.Range("D25").Formula = '="78.88"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +4.35%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("E28").Value = "  +3.21%  "

$ws.Range("D29").Formula = '="2.01"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +7.09%  "

$ws.Range("D30").Formula = '="7.47"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +7.08%  "

$ws.Range("D31").Formula = '="25.49"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +4.71%  "

$ws.Range("D32").Formula = '="5.92"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("E33").Value = "  +5.95%  "

$ws.Range("D34").Formula = '="5.69"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +7.51%  "

$ws.Range("D35").Formula = '="0.943"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +8.04%  "

$ws.Range("D36").Formula = '="2.07"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +6.08%  "

$ws.Range("D37").Formula = '="0.0₃0700"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +15.70%  "

$ws.Range("D38").Formula = '="48.26"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").Formula = '="8.72"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +7.72%  "

$ws.Range("E40").Value = "  +14.27%  "

$ws.Range("D41").Formula = '="379.90"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +9.53%  "

$ws.Range("E42").Value = "  +4.97%  "

$ws.Range("D44").Formula = '="2.698.38"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +5.08%  "

$ws.Range("D46").Formula = '="124.07"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +6.51%  "

$ws.Range("E47").Value = "  +5.43%  "

$ws.Range("E48").Value = "  +4.09%  "

$ws.Range("E49").Value = "  +2.93%  "

$ws.Range("E50").Value = "  +3.22%  "

$ws.Range("D51").Formula = '="2.00"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +5.36%  "
